# Updated filtering criteria to filter on the basis of a single keyword only.
# This re-runs the job-listing filter against the source data, which now
# yields a different (and smaller) set of matching rows. We clear out the
# old filtered results (rows 2-10, columns A:K) and write in the new
# filtered result set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous result block (rows 2-10 used columns A,B,C,D,E,G,H,I,J,K).
# Rows 11-13 (CSL Plasma / System Services / Indiana Blood Phlebotomy) are
# untouched by the new single-keyword filter, so they are left exactly as-is.
$ws.Range("A2:K10").ClearContents()

# New filtered rows: just a Title, a short trailing snippet of Experience,
# and the constant Employment_Type / Seniority_Level values.
$ws.Cells.Item(2, 1).Value = "Grand Junction Elementary School Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(2, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(2, 8).Value = "Full-time"
$ws.Cells.Item(2, 9).Value = "Entry level"

$ws.Cells.Item(3, 1).Value = "Grand Junction Elementary School Reading Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(3, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(3, 8).Value = "Full-time"
$ws.Cells.Item(3, 9).Value = "Entry level"

$ws.Cells.Item(4, 1).Value = "Grand Junction Elementary School Science Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(4, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(4, 8).Value = "Full-time"
$ws.Cells.Item(4, 9).Value = "Entry level"

$ws.Cells.Item(5, 1).Value = "Grand Junction Elementary School Writing Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(5, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(5, 8).Value = "Full-time"
$ws.Cells.Item(5, 9).Value = "Entry level"

$ws.Cells.Item(6, 1).Value = "Grand Junction Elementary School Math Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(6, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(6, 8).Value = "Full-time"
$ws.Cells.Item(6, 9).Value = "Entry level"

$ws.Cells.Item(7, 1).Value = "Grand Junction Middle School Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(7, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(7, 8).Value = "Full-time"
$ws.Cells.Item(7, 9).Value = "Entry level"

$ws.Cells.Item(8, 1).Value = "Grand Junction Middle School Reading Comprehension Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(8, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(8, 8).Value = "Full-time"
$ws.Cells.Item(8, 9).Value = "Entry level"

$ws.Cells.Item(9, 1).Value = "Grand Junction Middle School Reading Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(9, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(9, 8).Value = "Full-time"
$ws.Cells.Item(9, 9).Value = "Entry level"

$ws.Cells.Item(10, 1).Value = "Grand Junction High School Biology Certified Teacher Varsity Tutors, a Nerdy Company`r`n              `r`n          `r`n            `r`n              Grand Junction, CO"
$ws.Cells.Item(10, 5).Value = "experience. Show more Show less "
$ws.Cells.Item(10, 8).Value = "Full-time"
$ws.Cells.Item(10, 9).Value = "Entry level"
